# Sync non-localizable rule data: the "BannedPaths" rule row (row 35) moves
# down below the Oak-index rules (rows 36-40 shift up to 35-39), is renamed
# to the singular "BannedPath", its severity changes from Blocker to
# Critical, and it keeps having no Tags value (row 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the four "Oak index" rule rows up by one (36->35 ... 40->39).
$ws.Cells.Item(35,1).Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Cells.Item(35,2).Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Cells.Item(35,3).Value = "Bug"
$ws.Cells.Item(35,4).Value = "Blocker"
$ws.Cells.Item(35,5).Value = "aem,cloud-service-compatibility"

$ws.Cells.Item(36,1).Value = "IndexType"
$ws.Cells.Item(36,2).Value = "Custom Search Index Definition Nodes Must Use the Index Type lucene"
$ws.Cells.Item(36,3).Value = "Bug"
$ws.Cells.Item(36,4).Value = "Blocker"
$ws.Cells.Item(36,5).Value = "aem,cloud-service-compatibility"

$ws.Cells.Item(37,1).Value = "IndexAsyncProperty"
$ws.Cells.Item(37,2).Value = "Custom Lucene Oak Indexes must not be synchronous"
$ws.Cells.Item(37,3).Value = "Bug"
$ws.Cells.Item(37,4).Value = "Blocker"
$ws.Cells.Item(37,5).Value = "aem,cloud-service-compatibility"

$ws.Cells.Item(38,1).Value = "IndexTikaNode"
$ws.Cells.Item(38,2).Value = "Custom Oak indexes must have a tika configuration"
$ws.Cells.Item(38,3).Value = "Bug"
$ws.Cells.Item(38,4).Value = "Blocker"
$ws.Cells.Item(38,5).Value = "aem,cloud-service-compatibility"

$ws.Cells.Item(39,1).Value = "IndexDamAssetLucene"
$ws.Cells.Item(39,2).Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Cells.Item(39,3).Value = "Bug"
$ws.Cells.Item(39,4).Value = "Blocker"
$ws.Cells.Item(39,5).Value = "aem,cloud-service-compatibility"

# Former row 35 (BannedPaths) now lands at row 40, renamed singular, with a
# bumped severity and no Tags entry.
$ws.Cells.Item(40,1).Value = "BannedPath"
$ws.Cells.Item(40,2).Value = "Customer packages should not install content under /libs"
$ws.Cells.Item(40,3).Value = "Bug"
$ws.Cells.Item(40,4).Value = "Critical"
$ws.Cells.Item(40,5).ClearContents()

# Update the saved selection/cursor position to A37.
$ws.Range("A37").Select()
